# Weekly data refresh: insert two new daily price records at the top of the
# data block (rows 293-294) and push the existing records down by two rows
# (old 293-388 -> new 295-390), growing the sheet from A1:R388 to A1:R390.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 293:388 down by two rows, creating two blank rows
# at 293:294 for the new records (this also pushes the old R388 row to R390).
$ws.Rows("293:294").Insert()

# New row 293
$ws.Range("A293").Value = 8
$ws.Range("B293").Value = "Terminal La Palmera de La Serena"
$ws.Range("C293").Value = "Coquimbo"
$ws.Range("D293").Value = 44809
$ws.Range("E293").Value = 4
$ws.Range("F293").Value = 100112032
$ws.Range("G293").Value = "Zapallo italiano"
$ws.Range("H293").Value = "Bola 8"
$ws.Range("I293").Value = "Primera"
$ws.Range("J293").Value = 600
$ws.Range("K293").Value = 18000
$ws.Range("L293").Value = 19000
$ws.Range("M293").Value = 18500
$ws.Range("N293").Value = "`$/caja 50 unidades"
$ws.Range("O293").Value = "Región de Arica y Parinacota"
$ws.Range("P293").Value = 370
$ws.Range("Q293").Value = 50
$ws.Range("R293").Value = "Hortaliza"

# New row 294
$ws.Range("A294").Value = 8
$ws.Range("B294").Value = "Terminal La Palmera de La Serena"
$ws.Range("C294").Value = "Coquimbo"
$ws.Range("D294").Value = 44809
$ws.Range("E294").Value = 4
$ws.Range("F294").Value = 100112032
$ws.Range("G294").Value = "Zapallo italiano"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 600
$ws.Range("K294").Value = 20000
$ws.Range("L294").Value = 21000
$ws.Range("M294").Value = 20500
$ws.Range("N294").Value = "`$/caja 50 unidades"
$ws.Range("O294").Value = "Región de Arica y Parinacota"
$ws.Range("P294").Value = 410
$ws.Range("Q294").Value = 50
$ws.Range("R294").Value = "Hortaliza"
